$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.584.28'
$ws.Range('E2').Value = '  +5.66%  '
$ws.Range('D3').Value = '3.183.88'
$ws.Range('E3').Value = '  +3.06%  '
$ws.Range('D5').Value = "'401.26"
$ws.Range('E5').Value = '  +3.14%  '
$ws.Range('D6').Value = "'108.78"
$ws.Range('E6').Value = '  +5.40%  '
$ws.Range('E7').Value = '  +1.27%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = "'0.618"
$ws.Range('E9').Value = '  +5.15%  '
$ws.Range('D10').Value = "'39.04"
$ws.Range('E10').Value = '  +5.44%  '
$ws.Range('E11').Value = '  +1.59%  '
$ws.Range('D12').Value = "'0.0883"
$ws.Range('E12').Value = '  +2.68%  '
$ws.Range('D13').Value = '3.673.62'
$ws.Range('E13').Value = '  +2.60%  '
$ws.Range('D14').Value = "'19.03"
$ws.Range('E14').Value = '  +1.96%  '
$ws.Range('D15').Value = "'8.05"
$ws.Range('E15').Value = '  +3.11%  '
$ws.Range('E16').Value = '  +8.84%  '
$ws.Range('D17').Value = '3.184.03'
$ws.Range('E17').Value = '  +2.81%  '
$ws.Range('D18').Value = "'10.51"
$ws.Range('E18').Value = '  -1.22%  '
$ws.Range('D19').Value = '54.488.75'
$ws.Range('E19').Value = '  +5.16%  '
$ws.Range('D20').Value = "'3.33"
$ws.Range('E20').Value = '  +4.24%  '
$ws.Range('D21').Value = "'12.90"
$ws.Range('E21').Value = '  +3.40%  '
$ws.Range('D22').Value = '0.0₃0994'
$ws.Range('E22').Value = '  +2.56%  '
$ws.Range('D23').Value = "'72.24"
$ws.Range('E23').Value = '  +3.18%  '
$ws.Range('D24').Value = "'274.83"
$ws.Range('E24').Value = '  +2.29%  '
$ws.Range('E25').Value = '  +4.10%  '
$ws.Range('E26').Value = '  -1.96%  '
$ws.Range('D27').Value = "'27.74"
$ws.Range('E27').Value = '  +2.44%  '
$ws.Range('D28').Value = "'7.44"
$ws.Range('E28').Value = '  +3.27%  '
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  +4.56%  '
$ws.Range('E32').Value = '  +6.60%  '
$ws.Range('D33').Value = "'0.0510"
$ws.Range('E33').Value = '  +13.54%  '
$ws.Range('E34').Value = '  +4.37%  '
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('D36').Value = "'50.87"
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('D37').Value = "'3.66"
$ws.Range('E37').Value = '  +7.49%  '
$ws.Range('D38').Value = "'0.999"
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('D39').Value = "'2.84"
$ws.Range('E39').Value = '  +10.34%  '
$ws.Range('D40').Value = "'4.18"
$ws.Range('E40').Value = '  +12.85%  '
$ws.Range('E41').Value = '  +3.27%  '
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('D43').Value = "'17.35"
$ws.Range('E43').Value = '  +1.90%  '
$ws.Range('D44').Value = "'130.58"
$ws.Range('E44').Value = '  +2.66%  '
$ws.Range('D45').Value = "'0.118"
$ws.Range('E45').Value = '  +1.42%  '
$ws.Range('D46').Value = "'22.35"
$ws.Range('E46').Value = '  +1.03%  '
$ws.Range('D47').Value = "'2.46"
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('D48').Value = "'2.06"
$ws.Range('E48').Value = '  -0.89%  '
$ws.Range('D49').Value = '2.091.83'
$ws.Range('E49').Value = '  +2.28%  '
$ws.Range('D50').Value = "'0.0347"
$ws.Range('E50').Value = '  +8.86%  '
$ws.Range('E51').Value = '  +11.58%  '
